$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells F1:H1, copying formatting from the existing header (A1) ---
$ws.Range("A1").Copy($ws.Range("F1:H1"))
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# --- Boolean "outlier" flags for each algorithm column, rows 2-18 ---
$knnOutliers = @(0,0,1,0,0,0,0,1,0,0,0,0,0,0,0,0,0)
$svmOutliers = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
$rfOutliers  = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt 17; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = [bool]$knnOutliers[$i]
    $ws.Cells.Item($row, 7).Value = [bool]$svmOutliers[$i]
    $ws.Cells.Item($row, 8).Value = [bool]$rfOutliers[$i]
}
